$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold numeric-looking text (e.g. "1.040",
# "27.755.21", "  +2.62%  ") that must stay plain text, not be coerced to
# numbers by Excel. Force Text format before assigning those values.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.745.33'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.54%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.864.02'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.040'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +2.96%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.69'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.49%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.91%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4425'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +2.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3795'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07469'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8836'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.64'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.880.38'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -14.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.556'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.56%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.755'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07240'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.88'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +3.14%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.043'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.95%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009121'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +2.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.037'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.56'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.760.97'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.46%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.321'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.74%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.33'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.77%  '
$ws.Range('B24').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C24').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.100.54'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -12.39%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.008'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +6.48%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '159.07'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +3.07%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.83'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.983'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +4.10%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.316'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.96'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +2.84%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09062'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7758'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.61%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.209'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.018'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.65%  '
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.564'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +2.75%  '
$ws.Range('B36').Value = 'Frax'
$ws.Range('C36').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.038'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.98%  '
$ws.Range('B37').Value = 'TrustWalletToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.148'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.16%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01990'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +3.21%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05335'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.01%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.859'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +4.02%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.5197'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +1.43%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1692'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.32%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.859'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +5.76%  '
$ws.Range('B44').Value = 'Aptos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.652'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +3.31%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '110.10'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.72%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.69'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.06572'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +5.53%  '
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.720'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.26%  '
$ws.Range('B49').Value = 'Decentraland'
$ws.Range('C49').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4706'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.54%  '
$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.908'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.37%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.79'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.58%  '
